$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B/C/E: plain text cells, safe to assign directly. ---
$ws.Range('E2').Value = '  -5.02%  '
$ws.Range('E3').Value = '  -6.67%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('E5').Value = '  -3.13%  '
$ws.Range('E6').Value = '  -7.62%  '
$ws.Range('E8').Value = '  -6.77%  '
$ws.Range('E9').Value = '  -2.76%  '
$ws.Range('E10').Value = '  -6.21%  '
$ws.Range('E11').Value = '  -3.15%  '
$ws.Range('E12').Value = '  -3.20%  '
$ws.Range('E13').Value = '  -6.63%  '
$ws.Range('E14').Value = '  -6.10%  '
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('E16').Value = '  -6.88%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('E17').Value = '  -5.11%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('E18').Value = '  -6.98%  '
$ws.Range('E19').Value = '  -5.81%  '
$ws.Range('E20').Value = '  -6.31%  '
$ws.Range('E21').Value = '  -6.71%  '
$ws.Range('E22').Value = '  -5.48%  '
$ws.Range('E23').Value = '  -8.42%  '
$ws.Range('E24').Value = '  -4.71%  '
$ws.Range('E25').Value = '  -4.30%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('E28').Value = '  -5.94%  '
$ws.Range('E29').Value = '  -5.76%  '
$ws.Range('E30').Value = '  -7.75%  '
$ws.Range('E31').Value = '  -9.20%  '
$ws.Range('E32').Value = '  -7.15%  '
$ws.Range('E33').Value = '  -9.29%  '
$ws.Range('E34').Value = '  -8.24%  '
$ws.Range('E35').Value = '  -8.24%  '
$ws.Range('E36').Value = '  -4.01%  '
$ws.Range('E37').Value = '  -3.49%  '
$ws.Range('E38').Value = '  -7.89%  '
$ws.Range('E39').Value = '  -1.82%  '
$ws.Range('E40').Value = '  -8.14%  '
$ws.Range('E41').Value = '  -2.72%  '
$ws.Range('E42').Value = '  -6.40%  '
$ws.Range('E43').Value = '  -7.21%  '
$ws.Range('E44').Value = '  -6.72%  '
$ws.Range('E46').Value = '  -6.52%  '
$ws.Range('E48').Value = '  -6.49%  '
$ws.Range('E49').Value = '  -4.38%  '
$ws.Range('E50').Value = '  -7.25%  '
$ws.Range('E51').Value = '  -11.02%  '

# --- Column D ("Price"): values are dot-grouped, numeric-looking
# strings (e.g. "573.24", "59.976.17"). A direct .Value assignment
# lets Excel auto-coerce simple decimals into real numbers, which
# both changes the cell type and can silently drop significant
# trailing zeros (e.g. "0.660" -> 0.66). To keep these as literal
# text (matching the source workbook's inlineStr cells) without
# leaving a left-over number-format style on the cell, route each
# value through a text formula and then flatten that formula down
# to a plain value in place via Copy + PasteSpecial (values only).
$priceUpdates = @(
    @{ Cell = 'D2'; Text = '59.976.17' },
    @{ Cell = 'D3'; Text = '2.969.14' },
    @{ Cell = 'D5'; Text = '573.24' },
    @{ Cell = 'D6'; Text = '124.91' },
    @{ Cell = 'D8'; Text = '2.963.67' },
    @{ Cell = 'D9'; Text = '0.501' },
    @{ Cell = 'D12'; Text = '0.439' },
    @{ Cell = 'D14'; Text = '32.49' },
    @{ Cell = 'D16'; Text = '3.452.41' },
    @{ Cell = 'D17'; Text = '59.847.83' },
    @{ Cell = 'D18'; Text = '2.959.45' },
    @{ Cell = 'D19'; Text = '6.18' },
    @{ Cell = 'D20'; Text = '432.41' },
    @{ Cell = 'D21'; Text = '13.08' },
    @{ Cell = 'D22'; Text = '0.660' },
    @{ Cell = 'D23'; Text = '6.98' },
    @{ Cell = 'D24'; Text = '12.74' },
    @{ Cell = 'D25'; Text = '79.02' },
    @{ Cell = 'D27'; Text = '0.999' },
    @{ Cell = 'D28'; Text = '2.52' },
    @{ Cell = 'D29'; Text = '7.22' },
    @{ Cell = 'D32'; Text = '25.26' },
    @{ Cell = 'D33'; Text = '0.0929' },
    @{ Cell = 'D34'; Text = '2.18' },
    @{ Cell = 'D35'; Text = '0.949' },
    @{ Cell = 'D36'; Text = '5.58' },
    @{ Cell = 'D37'; Text = '49.57' },
    @{ Cell = 'D38'; Text = '0.0₃0653' },
    @{ Cell = 'D39'; Text = '7.94' },
    @{ Cell = 'D40'; Text = '0.0357' },
    @{ Cell = 'D42'; Text = '380.51' },
    @{ Cell = 'D44'; Text = '2.626.07' },
    @{ Cell = 'D47'; Text = '118.41' },
    @{ Cell = 'D50'; Text = '23.33' },
    @{ Cell = 'D51'; Text = '31.17' }
)

foreach ($u in $priceUpdates) {
    $rng = $ws.Range($u.Cell)
    $rng.Formula = '="' + $u.Text + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = 0
